$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Component" column (F) cells that referenced the now
#        removed "AmountWithCurrencyComponent" entry so they use the plain
#        "Currency" component instead. ---
foreach ($r in @(12, 14, 16, 19, 31, 33, 35, 38, 50, 52, 54, 57)) {
    $ws.Range("F$r").Value = "Currency"
}

# --- 2. Mark the "Document-Support" column (I) with "Extended" for every
#        row that represents a normal (non Custom-Component / non
#        Report-Preupload) field. Rows 7, 9, 17, 26, 36, 45, 55 and 64 are
#        intentionally skipped (Custom*/Report Preupload components). ---
$ws.Range("I2:I6").Value = "Extended"
$ws.Range("I8").Value = "Extended"
$ws.Range("I11:I16").Value = "Extended"
$ws.Range("I18:I25").Value = "Extended"
$ws.Range("I27:I28").Value = "Extended"
$ws.Range("I30:I35").Value = "Extended"
$ws.Range("I37:I44").Value = "Extended"
$ws.Range("I46:I47").Value = "Extended"
$ws.Range("I49:I54").Value = "Extended"
$ws.Range("I56:I57").Value = "Extended"

# Rows 58-66 currently use a slightly different (darker) fill style (s=19)
# for columns H-L. Column I there needs to switch to the normal style
# (s=12, the same one used by I2:I57) -- copy that formatting across before
# writing the values so the existing style is reused instead of a new one
# being created.
$ws.Range("I2").Copy()
$ws.Range("I58:I66").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I58:I63").Value = "Extended"
$ws.Range("I65:I66").Value = "Extended"
# I64 corresponds to a "Custom EuTaxonomyAlignedActivitiesComponent" row and
# stays blank (only its style changed above).

# --- 3. Turn the data range into an Excel Table-like AutoFilter and record
#        the associated hidden _FilterDatabase defined name, like Excel
#        does when the user applies Data > Filter. ---
$dataRange = $ws.Range("A1:L66")
$dataRange.AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $dataRange)
$filterName.Visible = $false
